$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D, E
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.663.37'
$ws.Range("E2").Value = '  -3.04%  '

# Row 3: update D, E
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.312.53'
$ws.Range("E3").Value = '  -4.12%  '

# Row 4: update E
$ws.Range("E4").Value = '  -0.02%  '

# Row 5: update D, E
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.32'
$ws.Range("E5").Value = '  -2.92%  '

# Row 6: update D, E
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.70'
$ws.Range("E6").Value = '  -7.83%  '

# Row 7: update D, E
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.516'
$ws.Range("E7").Value = '  -3.56%  '

# Row 8: update E
$ws.Range("E8").Value = '  -0.01%  '

# Row 9: update D, E
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.471'
$ws.Range("E9").Value = '  -5.03%  '

# Row 10: update D, E
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0786'
$ws.Range("E10").Value = '  -5.56%  '

# Row 11: update D, E
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.69'
$ws.Range("E11").Value = '  -9.80%  '

# Row 12: update E
$ws.Range("E12").Value = '  -0.35%  '

# Row 13: update D, E
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.674.09'
$ws.Range("E13").Value = '  -3.95%  '

# Row 14: update D, E
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.24'
$ws.Range("E14").Value = '  -6.52%  '

# Row 15: update D, E
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.41'
$ws.Range("E15").Value = '  -7.96%  '

# Row 16: update D, E
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.322.93'
$ws.Range("E16").Value = '  -3.18%  '

# Row 17: update D, E
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.738'
$ws.Range("E17").Value = '  -4.43%  '

# Row 18: update D, E
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.608.50'
$ws.Range("E18").Value = '  -3.06%  '

# Row 19: update D, E
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0881'
$ws.Range("E19").Value = '  -4.37%  '

# Row 20: update D, E
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.92'
$ws.Range("E20").Value = '  -5.15%  '

# Row 21: update D, E
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.29'
$ws.Range("E21").Value = '  -6.29%  '

# Row 22: update D, E
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.23'
$ws.Range("E22").Value = '  -6.33%  '

# Row 23: update D, E
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.95'
$ws.Range("E23").Value = '  -0.78%  '

# Row 24: update D, E
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.49'
$ws.Range("E24").Value = '  -6.83%  '

# Row 25: update E
$ws.Range("E25").Value = '  -0.07%  '

# Row 26: update D, E
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.78'
$ws.Range("E26").Value = '  -4.24%  '

# Row 27: update D, E
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.85'
$ws.Range("E27").Value = '  -4.84%  '

# Row 28: update D, E
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.18'
$ws.Range("E28").Value = '  -2.24%  '

# Row 29: update D, E
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.04'
$ws.Range("E29").Value = '  -5.11%  '

# Row 30: update D, E
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.86'
$ws.Range("E30").Value = '  -4.06%  '

# Row 31: update D, E
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '151.19'
$ws.Range("E31").Value = '  -3.86%  '

# Row 32: update E
$ws.Range("E32").Value = '  -0.15%  '

# Row 33: update D, E
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.95'
$ws.Range("E33").Value = '  -5.61%  '

# Row 34: update D, E
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.43'
$ws.Range("E34").Value = '  -0.42%  '

# Row 35: update D, E
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0701'
$ws.Range("E35").Value = '  -5.43%  '

# Row 36: update D, E
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.112'
$ws.Range("E36").Value = '  -1.78%  '

# Row 37: update D, E
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.70'
$ws.Range("E37").Value = '  -6.77%  '

# Row 38: update D, E
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0965'
$ws.Range("E38").Value = '  -3.56%  '

# Row 39: update D, E
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.11'
$ws.Range("E39").Value = '  -8.54%  '

# Row 40: update D, E
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.64'
$ws.Range("E40").Value = '  -7.47%  '

# Row 41: update D, E
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.68'
$ws.Range("E41").Value = '  -4.58%  '

# Row 42: update B, C, D, E
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.26'
$ws.Range("E42").Value = '  -2.80%  '

# Row 43: update B, C, D, E
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.958.58'
$ws.Range("E43").Value = '  -1.65%  '

# Row 44: update D, E
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0258'
$ws.Range("E44").Value = '  -5.62%  '

# Row 45: update D, E
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.71'
$ws.Range("E45").Value = '  -10.67%  '

# Row 46: update D, E
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.21'
$ws.Range("E46").Value = '  -1.67%  '

# Row 47: update D, E
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.60'
$ws.Range("E47").Value = '  -9.47%  '

# Row 48: update D, E
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.533.71'
$ws.Range("E48").Value = '  -4.40%  '

# Row 49: update D, E
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '90.57'
$ws.Range("E49").Value = '  -3.66%  '

# Row 50: update D, E
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.31'
$ws.Range("E50").Value = '  -6.76%  '

# Row 51: update D, E
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '48.49'
$ws.Range("E51").Value = '  -5.93%  '

Write-Host "Applied cryptos update"